# Horarios actualizados Linea 141 - 375
$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 08:19:34"
$ws1.Range("A3").Value = "Total filas: 96"

$ws1.Cells.Item(77, 1).Value = "07:58:19"
$ws1.Cells.Item(77, 2).Value = "08:42"
$ws1.Cells.Item(77, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(77, 4).Value = 44
$ws1.Cells.Item(77, 5).Value = "LP1912"

$ws1.Cells.Item(78, 1).Value = "06:57:11"
$ws1.Cells.Item(78, 2).Value = "08:42"
$ws1.Cells.Item(78, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(78, 4).Value = 105
$ws1.Cells.Item(78, 5).Value = "LP1912"

$ws1.Cells.Item(79, 1).Value = "07:19:37"
$ws1.Cells.Item(79, 2).Value = "08:43"
$ws1.Cells.Item(79, 3).Value = "14_ABASTO"
$ws1.Cells.Item(79, 4).Value = 84
$ws1.Cells.Item(79, 5).Value = "LP1912"

$ws1.Cells.Item(80, 1).Value = "08:19:33"
$ws1.Cells.Item(80, 2).Value = "08:43"
$ws1.Cells.Item(80, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(80, 4).Value = 24
$ws1.Cells.Item(80, 5).Value = "LP1912"

$ws1.Cells.Item(81, 1).Value = "08:19:33"
$ws1.Cells.Item(81, 2).Value = "08:53"
$ws1.Cells.Item(81, 3).Value = "10_OLMOS"
$ws1.Cells.Item(81, 4).Value = 34
$ws1.Cells.Item(81, 5).Value = "LP1912"

$ws1.Cells.Item(82, 1).Value = "06:57:11"
$ws1.Cells.Item(82, 2).Value = "08:54"
$ws1.Cells.Item(82, 3).Value = "17_ROMERO"
$ws1.Cells.Item(82, 4).Value = 117
$ws1.Cells.Item(82, 5).Value = "LP1912"

$ws1.Cells.Item(83, 1).Value = "07:19:37"
$ws1.Cells.Item(83, 2).Value = "09:01"
$ws1.Cells.Item(83, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(83, 4).Value = 102
$ws1.Cells.Item(83, 5).Value = "LP1912"

$ws1.Cells.Item(84, 1).Value = "08:19:33"
$ws1.Cells.Item(84, 2).Value = "09:02"
$ws1.Cells.Item(84, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(84, 4).Value = 43
$ws1.Cells.Item(84, 5).Value = "LP1912"

$ws1.Cells.Item(85, 1).Value = "07:45:49"
$ws1.Cells.Item(85, 2).Value = "09:03"
$ws1.Cells.Item(85, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(85, 4).Value = 78
$ws1.Cells.Item(85, 5).Value = "LP1912"

$ws1.Cells.Item(86, 1).Value = "07:19:37"
$ws1.Cells.Item(86, 2).Value = "09:10"
$ws1.Cells.Item(86, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(86, 4).Value = 111
$ws1.Cells.Item(86, 5).Value = "LP1912"

$ws1.Cells.Item(87, 1).Value = "07:19:37"
$ws1.Cells.Item(87, 2).Value = "09:16"
$ws1.Cells.Item(87, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(87, 4).Value = 117
$ws1.Cells.Item(87, 5).Value = "LP1912"

$ws1.Cells.Item(88, 1).Value = "07:58:19"
$ws1.Cells.Item(88, 2).Value = "09:17"
$ws1.Cells.Item(88, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(88, 4).Value = 79
$ws1.Cells.Item(88, 5).Value = "LP1912"

$ws1.Cells.Item(89, 1).Value = "07:45:49"
$ws1.Cells.Item(89, 2).Value = "09:21"
$ws1.Cells.Item(89, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(89, 4).Value = 96
$ws1.Cells.Item(89, 5).Value = "LP1912"

$ws1.Cells.Item(90, 1).Value = "07:45:49"
$ws1.Cells.Item(90, 2).Value = "09:22"
$ws1.Cells.Item(90, 3).Value = "17_ROMERO"
$ws1.Cells.Item(90, 4).Value = 97
$ws1.Cells.Item(90, 5).Value = "LP1912"

$ws1.Cells.Item(91, 1).Value = "07:45:49"
$ws1.Cells.Item(91, 2).Value = "09:23"
$ws1.Cells.Item(91, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(91, 4).Value = 98
$ws1.Cells.Item(91, 5).Value = "LP1912"

$ws1.Cells.Item(92, 1).Value = "07:58:19"
$ws1.Cells.Item(92, 2).Value = "09:23"
$ws1.Cells.Item(92, 3).Value = "17_ROMERO"
$ws1.Cells.Item(92, 4).Value = 85
$ws1.Cells.Item(92, 5).Value = "LP1912"

$ws1.Cells.Item(93, 1).Value = "08:19:33"
$ws1.Cells.Item(93, 2).Value = "09:25"
$ws1.Cells.Item(93, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(93, 4).Value = 66
$ws1.Cells.Item(93, 5).Value = "LP1912"

$ws1.Cells.Item(94, 1).Value = "07:45:49"
$ws1.Cells.Item(94, 2).Value = "09:32"
$ws1.Cells.Item(94, 3).Value = "15_ABASTO"
$ws1.Cells.Item(94, 4).Value = 107
$ws1.Cells.Item(94, 5).Value = "LP1912"

$ws1.Cells.Item(95, 1).Value = "07:45:49"
$ws1.Cells.Item(95, 2).Value = "09:33"
$ws1.Cells.Item(95, 3).Value = "10_OLMOS"
$ws1.Cells.Item(95, 4).Value = 108
$ws1.Cells.Item(95, 5).Value = "LP1912"

$ws1.Cells.Item(96, 1).Value = "07:45:49"
$ws1.Cells.Item(96, 2).Value = "09:41"
$ws1.Cells.Item(96, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(96, 4).Value = 116
$ws1.Cells.Item(96, 5).Value = "LP1912"

$ws1.Cells.Item(97, 1).Value = "07:58:19"
$ws1.Cells.Item(97, 2).Value = "09:42"
$ws1.Cells.Item(97, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(97, 4).Value = 104
$ws1.Cells.Item(97, 5).Value = "LP1912"

$ws1.Cells.Item(98, 1).Value = "07:58:19"
$ws1.Cells.Item(98, 2).Value = "09:43"
$ws1.Cells.Item(98, 3).Value = "14_ABASTO"
$ws1.Cells.Item(98, 4).Value = 105
$ws1.Cells.Item(98, 5).Value = "LP1912"

$ws1.Cells.Item(99, 1).Value = "07:58:19"
$ws1.Cells.Item(99, 2).Value = "09:52"
$ws1.Cells.Item(99, 3).Value = "15_ABASTO"
$ws1.Cells.Item(99, 4).Value = 114
$ws1.Cells.Item(99, 5).Value = "LP1912"

$ws1.Cells.Item(100, 1).Value = "08:19:33"
$ws1.Cells.Item(100, 2).Value = "10:10"
$ws1.Cells.Item(100, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(100, 4).Value = 111
$ws1.Cells.Item(100, 5).Value = "LP1912"

$ws1.Cells.Item(101, 1).Value = "08:19:33"
$ws1.Cells.Item(101, 2).Value = "10:12"
$ws1.Cells.Item(101, 3).Value = "15_ABASTO"
$ws1.Cells.Item(101, 4).Value = 113
$ws1.Cells.Item(101, 5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 08:19:34"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 08:19:34"
$ws3.Range("A3").Value = "Total filas: 21"

$ws3.Cells.Item(23, 1).Value = "08:19:33"
$ws3.Cells.Item(23, 2).Value = "08:43"
$ws3.Cells.Item(23, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(23, 4).Value = 24
$ws3.Cells.Item(23, 5).Value = "L6173"

$ws3.Cells.Item(24, 1).Value = "07:19:37"
$ws3.Cells.Item(24, 2).Value = "09:08"
$ws3.Cells.Item(24, 3).Value = "215D_LA PLATA"
$ws3.Cells.Item(24, 4).Value = 109
$ws3.Cells.Item(24, 5).Value = "L6203"

$ws3.Cells.Item(25, 1).Value = "07:58:19"
$ws3.Cells.Item(25, 2).Value = "09:09"
$ws3.Cells.Item(25, 3).Value = "215D_LA PLATA"
$ws3.Cells.Item(25, 4).Value = 71
$ws3.Cells.Item(25, 5).Value = "L6203"

$ws3.Cells.Item(26, 1).Value = "08:19:33"
$ws3.Cells.Item(26, 2).Value = "10:03"
$ws3.Cells.Item(26, 3).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(26, 4).Value = 104
$ws3.Cells.Item(26, 5).Value = "L6173"

